$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (columns D, L, M, N, O, P, Q, R, S, T)
# derived from the permutation of the original rows described by the diff.

$data = @{
    2  = @{ D=44358; L="Primera"; M=100; N=11000; O=12000; P=11500; Q="`$/caja 18 kilos granel";      R="Región de O'Higgins"; S=639; T=18 }
    3  = @{ D=44698; L="Primera"; M=50;  N=10000; O=10000; P=10000; Q="`$/caja 18 kilos granel";      R="Región de O'Higgins"; S=556; T=18 }
    4  = @{ D=44299; L="Primera"; M=100; N=10000; O=11000; P=10500; Q="`$/caja 18 kilos granel";      R="Región del Maule";    S=583; T=18 }
    5  = @{ D=44299; L="Segunda"; M=50;  N=9000;  O=9000;  P=9000;  Q="`$/caja 18 kilos granel";      R="Región del Maule";    S=500; T=18 }
    6  = @{ D=44363; L="Primera"; M=100; N=9000;  O=10000; P=9500;  Q="`$/caja 15 kilos empedrada";   R="Región de O'Higgins"; S=633; T=15 }
    7  = @{ D=44776; L="Primera"; M=50;  N=10000; O=10000; P=10000; Q="`$/bandeja 18 kilos granel";   R="Región de O'Higgins"; S=556; T=18 }
    8  = @{ D=44776; L="Segunda"; M=50;  N=8000;  O=8000;  P=8000;  Q="`$/bandeja 18 kilos granel";   R="Región de O'Higgins"; S=444; T=18 }
    11 = @{ D=44425; L="Primera"; M=100; N=12000; O=13000; P=12500; Q="`$/bandeja 18 kilos granel";   R="Región de O'Higgins"; S=694; T=18 }
    12 = @{ D=44316; L="Primera"; M=100; N=9000;  O=10000; P=9500;  Q="`$/caja 18 kilos granel";      R="Región de O'Higgins"; S=528; T=18 }
    13 = @{ D=44272; L="Primera"; M=100; N=9000;  O=10000; P=9500;  Q="`$/caja 15 kilos granel";      R="Región de O'Higgins"; S=633; T=15 }
    14 = @{ D=44272; L="Segunda"; M=50;  N=8000;  O=8000;  P=8000;  Q="`$/caja 15 kilos granel";      R="Región de O'Higgins"; S=533; T=15 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
